$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.745.24"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "'3.133.20"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'529.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "'138.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'3.131.08"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "'0.448"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("D10").Value = "'7.21"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "'0.397"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("D13").Value = "'3.673.38"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "'25.60"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "'57.865.96"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'3.136.73"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "'12.76"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "'7.94"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "'352.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.88%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'68.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "'0.169"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "'7.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.03%  "
$ws.Range("D31").Value = "'6.20"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.98%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").Value = "'21.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").Value = "'1.18"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("E35").Value = "  +7.54%  "
$ws.Range("D36").Value = "'158.92"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "'6.16"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").Value = "'26.48"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +6.65%  "
$ws.Range("E42").Value = "  +6.73%  "
$ws.Range("D43").Value = "'0.702"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").Value = "'3.173.74"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'36.59"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0270"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.03%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'2.313.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "'0.963"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "'20.36"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.97%  "
